$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2135627530364372
$ws.Range("C2").Value = 0.5232793522267206
$ws.Range("J2").Value = 0.01214574898785425
$ws.Range("P2").Value = 0.1538461538461539
$ws.Range("S2").Value = 0.09716599190283401

# Row 3
$ws.Range("B3").Value = 0.02407407407407407
$ws.Range("C3").Value = 0.04074074074074074
$ws.Range("J3").Value = 0.03148148148148148
$ws.Range("P3").Value = 0.687037037037037
$ws.Range("S3").Value = 0.2166666666666667

# Row 4
$ws.Range("J4").Value = 0.032
$ws.Range("P4").Value = 0.656
$ws.Range("S4").Value = 0.312

# Row 6
$ws.Range("B6").Value = 0.09838472834067548
$ws.Range("D6").Value = 0.01615271659324523
$ws.Range("F6").Value = 0.08076358296622614
$ws.Range("J6").Value = 0.2452276064610866
$ws.Range("O6").Value = 0.01908957415565345
$ws.Range("Q6").Value = 0.1439060205580029
$ws.Range("R6").Value = 0.05286343612334802
$ws.Range("S6").Value = 0.3436123348017621

# Row 7
$ws.Range("B7").Value = 0.1127214170692432
$ws.Range("D7").Value = 0.01932367149758454
$ws.Range("E7").Value = 0.001610305958132045
$ws.Range("F7").Value = 0.0644122383252818
$ws.Range("J7").Value = 0.1384863123993559
$ws.Range("O7").Value = 0.02576489533011272
$ws.Range("Q7").Value = 0.14170692431562
$ws.Range("R7").Value = 0.07085346215780998
$ws.Range("S7").Value = 0.4251207729468599

# Row 8
$ws.Range("B8").Value = 0.08908839779005524
$ws.Range("D8").Value = 0.01450276243093923
$ws.Range("E8").Value = 0.002071823204419889
$ws.Range("F8").Value = 0.05524861878453038
$ws.Range("J8").Value = 0.1270718232044199
$ws.Range("O8").Value = 0.0138121546961326
$ws.Range("Q8").Value = 0.175414364640884
$ws.Range("R8").Value = 0.09875690607734806
$ws.Range("S8").Value = 0.4240331491712707

# Row 9
$ws.Range("B9").Value = 0.09777015437392796
$ws.Range("D9").Value = 0.02229845626072041
$ws.Range("F9").Value = 0.0686106346483705
$ws.Range("J9").Value = 0.1200686106346484
$ws.Range("O9").Value = 0.01886792452830189
$ws.Range("Q9").Value = 0.1749571183533448
$ws.Range("R9").Value = 0.08747855917667238
$ws.Range("S9").Value = 0.4099485420240137

# Row 10
$ws.Range("B10").Value = 0.1122742737503271
$ws.Range("D10").Value = 0.01910494634912327
$ws.Range("E10").Value = 0.001308557969118032
$ws.Range("F10").Value = 0.0698769955509029
$ws.Range("J10").Value = 0.1345197592253337
$ws.Range("O10").Value = 0.01544098403559278
$ws.Range("Q10").Value = 0.1884323475529966
$ws.Range("R10").Value = 0.07825176655325831
$ws.Range("S10").Value = 0.3807903690133473

# Row 11
$ws.Range("G11").Value = 0.1445396145610278
$ws.Range("J11").Value = 0.07815845824411134
$ws.Range("K11").Value = 0.1905781584582441
$ws.Range("L11").Value = 0.5674518201284796
$ws.Range("S11").Value = 0.01927194860813704

# Row 12
$ws.Range("G12").Value = 0.7423971377459749
$ws.Range("J12").Value = 0.1735241502683363
$ws.Range("K12").Value = 0.007155635062611807
$ws.Range("L12").Value = 0.03935599284436494
$ws.Range("S12").Value = 0.03756708407871199

# Row 13
$ws.Range("G13").Value = 0.7241379310344828
$ws.Range("J13").Value = 0.2068965517241379
$ws.Range("S13").Value = 0.06896551724137931

# Row 15
$ws.Range("F15").Value = 0.03047091412742382
$ws.Range("H15").Value = 0.1620498614958449
$ws.Range("I15").Value = 0.0817174515235457
$ws.Range("J15").Value = 0.3670360110803324
$ws.Range("K15").Value = 0.06509695290858726
$ws.Range("M15").Value = 0.009695290858725761
$ws.Range("N15").Value = 0.001385041551246537
$ws.Range("O15").Value = 0.08448753462603878
$ws.Range("S15").Value = 0.1980609418282548

# Row 16
$ws.Range("F16").Value = 0.01700680272108844
$ws.Range("H16").Value = 0.195578231292517
$ws.Range("I16").Value = 0.06462585034013606
$ws.Range("J16").Value = 0.391156462585034
$ws.Range("K16").Value = 0.1105442176870748
$ws.Range("M16").Value = 0.02040816326530612
$ws.Range("O16").Value = 0.07312925170068027
$ws.Range("S16").Value = 0.1275510204081633

# Row 17
$ws.Range("F17").Value = 0.01682692307692308
$ws.Range("H17").Value = 0.1931089743589744
$ws.Range("I17").Value = 0.08413461538461539
$ws.Range("J17").Value = 0.3958333333333333
$ws.Range("K17").Value = 0.1073717948717949
$ws.Range("M17").Value = 0.01682692307692308
$ws.Range("N17").Value = 0.001602564102564103
$ws.Range("O17").Value = 0.07131410256410256
$ws.Range("S17").Value = 0.1129807692307692

# Row 18
$ws.Range("F18").Value = 0.02083333333333333
$ws.Range("H18").Value = 0.1684027777777778
$ws.Range("I18").Value = 0.1059027777777778
$ws.Range("J18").Value = 0.4184027777777778
$ws.Range("K18").Value = 0.1024305555555556
$ws.Range("M18").Value = 0.01909722222222222
$ws.Range("O18").Value = 0.07465277777777778
$ws.Range("S18").Value = 0.09027777777777778

# Row 19
$ws.Range("F19").Value = 0.01835798062213156
$ws.Range("H19").Value = 0.2218255991840898
$ws.Range("I19").Value = 0.0815910249872514
$ws.Range("J19").Value = 0.3607853136155023
$ws.Range("K19").Value = 0.11116777154513
$ws.Range("M19").Value = 0.01708312085670576
$ws.Range("N19").Value = 0.001529831718510964
$ws.Range("O19").Value = 0.06986231514533402
$ws.Range("S19").Value = 0.1177970423253442
